# Corrected Excel import error handling
# Adds a new data row (row 5) to Sheet1 that duplicates the content of row 2
# (the "QnABot.001" sample question) but omits the "qid" value in column A,
# so the sheet can be used to exercise the import tool's handling of rows
# that are missing a qid. A hyperlink matching the one on I2 is added for I5,
# and the sheet view is updated to show the newly added row as selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate row 5 (B:M) with the same values as row 2, leaving A5 empty ---
$ws.Range("B5").Value = $ws.Range("B2").Value()
$ws.Range("C5").Value = $ws.Range("C2").Value()
$ws.Range("D5").Value = $ws.Range("D2").Value()
$ws.Range("E5").Value = $ws.Range("E2").Value()
$ws.Range("F5").Value = $ws.Range("F2").Value()
$ws.Range("G5").Value = $ws.Range("G2").Value()
$ws.Range("H5").Value = $ws.Range("H2").Value()
$ws.Range("I5").Value = $ws.Range("I2").Value()
$ws.Range("J5").Value = $ws.Range("J2").Value()
$ws.Range("K5").Value = $ws.Range("K2").Value()
$ws.Range("L5").Value = $ws.Range("L2").Value()
$ws.Range("M5").Value = $ws.Range("M2").Value()

# --- Match row 2's cell formatting (wrap-text styles, row height) ---
$ws.Range("J2").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("L2").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("M2").Copy()
$ws.Range("M5").PasteSpecial(-4122)

$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(2).RowHeight

# --- Re-create the image hyperlink on I5, matching I2/I3/I4, then restore I2's
#     cell format (Hyperlinks.Add applies its own built-in style by default) ---
$ws.Hyperlinks.Add($ws.Range("I5"), "https://images-na.ssl-images-amazon.com/images/I/61bze1WJhfL._AC_SL1024_.jpg")
$ws.Range("I2").Copy()
$ws.Range("I5").PasteSpecial(-4122)

# --- Update the sheet view to reflect scrolling down to show the new row ---
$ws.Range("A5").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
